$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume change (E) columns; rows 44/45 also swap coin name/link (B/C).

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = '29.661.81'
$dCell.ClearFormats()
$ws.Range("E2").Value = '  +0.55%  '

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = '1.616.69'
$dCell.ClearFormats()
$ws.Range("E3").Value = '  +0.65%  '

$ws.Range("E4").Value = '  -0.43%  '

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = '212.60'
$dCell.ClearFormats()
$ws.Range("E5").Value = '  +0.02%  '

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = '0.522'
$dCell.ClearFormats()
$ws.Range("E6").Value = '  +0.39%  '

$ws.Range("E7").Value = '  -0.47%  '

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = '28.87'
$dCell.ClearFormats()
$ws.Range("E8").Value = '  +8.61%  '

$ws.Range("E9").Value = '  +2.85%  '

$ws.Range("E10").Value = '  +1.70%  '

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0910'
$dCell.ClearFormats()
$ws.Range("E11").Value = '  -0.05%  '

$ws.Range("E12").Value = '  +0.76%  '

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = '1.611.91'
$dCell.ClearFormats()
$ws.Range("E13").Value = '  +0.26%  '

$ws.Range("E14").Value = '  +6.16%  '

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = '3.85'
$dCell.ClearFormats()
$ws.Range("E15").Value = '  +4.12%  '

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = '29.683.47'
$dCell.ClearFormats()
$ws.Range("E16").Value = '  +0.51%  '

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = '8.89'
$dCell.ClearFormats()
$ws.Range("E17").Value = '  +16.33%  '

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = '64.36'
$dCell.ClearFormats()
$ws.Range("E18").Value = '  +1.55%  '

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = '243.23'
$dCell.ClearFormats()
$ws.Range("E19").Value = '  +0.22%  '

$ws.Range("E20").Value = '  +2.72%  '

$ws.Range("E21").Value = '  -0.36%  '

$ws.Range("E22").Value = '  +3.10%  '

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = '9.68'
$dCell.ClearFormats()
$ws.Range("E23").Value = '  +5.80%  '

$ws.Range("E24").Value = '  +0.87%  '

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = '156.68'
$dCell.ClearFormats()
$ws.Range("E25").Value = '  +1.32%  '

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = '15.64'
$dCell.ClearFormats()
$ws.Range("E26").Value = '  +2.02%  '

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = '6.59'
$dCell.ClearFormats()
$ws.Range("E28").Value = '  +3.30%  '

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = '0.995'
$dCell.ClearFormats()
$ws.Range("E29").Value = '  -0.34%  '

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0484'
$dCell.ClearFormats()
$ws.Range("E30").Value = '  +2.28%  '

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = '3.31'
$dCell.ClearFormats()
$ws.Range("E31").Value = '  +3.05%  '

$ws.Range("E32").Value = '  +1.31%  '

$ws.Range("E33").Value = '  +3.40%  '

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = '1.434.93'
$dCell.ClearFormats()
$ws.Range("E34").Value = '  +1.28%  '

$ws.Range("E35").Value = '  +6.68%  '

$ws.Range("E36").Value = '  +1.71%  '

$ws.Range("E37").Value = '  +2.48%  '

$ws.Range("E38").Value = '  -0.80%  '

$ws.Range("E39").Value = '  +3.26%  '

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = '0.556'
$dCell.ClearFormats()
$ws.Range("E40").Value = '  +3.69%  '

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0504'
$dCell.ClearFormats()
$ws.Range("E41").Value = '  +4.73%  '

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = '1.97'
$dCell.ClearFormats()
$ws.Range("E43").Value = '  +0.50%  '

$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = '54.04'
$dCell.ClearFormats()
$ws.Range("E44").Value = '  +1.00%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = '69.48'
$dCell.ClearFormats()
$ws.Range("E45").Value = '  +5.86%  '

$ws.Range("E46").Value = '  -0.49%  '

$ws.Range("E47").Value = '  +19.14%  '

$ws.Range("E48").Value = '  +3.06%  '

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = '1.758.15'
$dCell.ClearFormats()
$ws.Range("E49").Value = '  +0.50%  '

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = '88.12'
$dCell.ClearFormats()
$ws.Range("E50").Value = '  +1.67%  '

$ws.Range("E51").Value = '  -0.86%  '
